$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 previously had redundant "value" header labels repeated across
# C1:F1 - clear those leftover cells.
$ws.Range("C1:F1").ClearContents()

# "Model" label becomes "production_function".
$ws.Range("A8").Value = "production_function"

# Insert a new parameter row right after it for the L-curve toggle.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").NumberFormat = "0.00E+00"
$ws.Range("B9").Value = 1

# The old "Deletion" row (now shifted down to row 17 by the insert above)
# is no longer needed - remove it entirely.
$ws.Rows.Item(17).Delete()

# Reflect the user's final action: switch to this sheet and select the
# last row in full.
$ws.Activate()
$ws.Rows.Item(17).Select()
